$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 289, shifting existing rows 289-386 down to 290-387.
$ws.Range("A289").EntireRow.Insert()

# Populate the new row 289 with the new observation data.
$ws.Range("A289").Value = 6
$ws.Range("B289").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C289").Value = "Metropolitana"
$ws.Range("D289").Value = 44524
$ws.Range("E289").Value = 13
$ws.Range("F289").Value = 100112044
$ws.Range("G289").Value = "Perejil"
$ws.Range("H289").Value = "Sin especificar"
$ws.Range("I289").Value = "Primera"
$ws.Range("J289").Value = 190
$ws.Range("K289").Value = 9500
$ws.Range("L289").Value = 10000
$ws.Range("M289").Value = 9684
$ws.Range("N289").Value = "`$/docena de atados"
$ws.Range("O289").Value = "Región Metropolitana"
$ws.Range("P289").Value = 3228
$ws.Range("Q289").Value = 3
$ws.Range("R289").Value = "Hortaliza"
